$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Module 09: Tissues — Keys to Success",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Module 09: Tissues and the Animal Body — Keys to Success", 2)

$d.Content.Find.Execute(
    "1. Mendelian Genetics Foundations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Homeostasis and Osmoregulation", 2)

$d.Content.Find.Execute(
    "Define gene, allele, genotype, and phenotype  Distinguish between dominant and recessive alleles  Explain homozygous and heterozygous conditions  Describe Mendel's experiments and his two laws   2. Law of Segregation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Define homeostasis and describe its role in maintaining a stable internal environment  Explain negative and positive feedback loops using biological examples  Define osmoregulation and explain how the body balances water and salt   2. Digestive System", 2)

$d.Content.Find.Execute(
    "Explain the law of segregation  Use Punnett squares to predict offspring genotypes and phenotypes  Calculate genotypic and phenotypic ratios for monohybrid crosses   3. Law of Independent Assortment",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe the primary function of the digestive system  Identify the main organs of the digestive tract and their specific roles (mouth, stomach, small intestine, large intestine)  Explain the difference between mechanical and chemical digestion  Describe the role of accessory organs (liver, pancreas, gallbladder)   3. Circulatory and Respiratory Systems", 2)

$d.Content.Find.Execute(
    "Explain the law of independent assortment  Perform dihybrid crosses using Punnett squares  Understand when independent assortment applies and its limitations   4. Extensions to Mendelian Genetics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Explain the function of the circulatory system in transporting nutrients and waste  Describe the pathway of blood through the human heart and lungs  Differentiate between arteries, veins, and capillaries  Explain how the respiratory system facilitates gas exchange (oxygen and carbon dioxide)   4. Endocrine System", 2)

$d.Content.Find.Execute(
    "Describe incomplete dominance and codominance  Explain multiple alleles using ABO blood types as an example  Understand polygenic inheritance and continuous variation  Explain pleiotropy (one gene affecting multiple traits)   5. Sex Linkage and Chromosomal Inheritance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Define hormones and describe how the endocrine system uses them to communicate  Identify major endocrine glands (pituitary, thyroid, adrenal, pancreas) and their primary functions  Explain how the endocrine system works with the nervous system to maintain homeostasis   5. Musculoskeletal System", 2)

$d.Content.Find.Execute(
    "Describe sex determination in humans  Explain sex-linked inheritance patterns  Predict outcomes of crosses involving X-linked traits  Understand why sex-linked disorders are more common in males   6. Pedigree Analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe the primary functions of the skeletal system (support, protection, movement, blood cell production)  Differentiate between the axial and appendicular skeleton  Differentiate between skeletal, smooth, and cardiac muscle tissue  Explain how muscles and bones work together to create movement   6. Nervous System", 2)

$d.Content.Find.Execute(
    "Interpret pedigree charts  Determine modes of inheritance from pedigrees  Identify carriers and affected individuals    Study Tips",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe the primary function of the nervous system in processing information  Differentiate between the central nervous system (CNS) and peripheral nervous system (PNS)  Identify the basic structure of a neuron and describe how signals are transmitted    Study Tips", 2)
